$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "Escolha das correções" checkbox/label groups ---
# Before:  I10(@) J10:L10(Free-Air) | N10(@) O10:R10(Bouguer Simples) | T10(@) U10:X10(Pressão Atmosférica)
# After:               L10(@) M10:P10(Free-air)   |   S10(@) T10:W10(Bouguer Simples)
# (third "Pressão Atmosférica" group removed entirely)

# Grab format (style) of a still-centered "@" cell (style index 2) before we touch it,
# so we can stamp that same format onto M10 (which currently carries the blank/style-11 format).
$ws.Range("N10").Copy()
$ws.Range("M10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Un-merge the three old groups, merge the two new groups
$ws.Range("J10:L10").UnMerge()
$ws.Range("O10:R10").UnMerge()
$ws.Range("U10:X10").UnMerge()
$ws.Range("M10:P10").Merge()
$ws.Range("T10:W10").Merge()

# Move/update the cell contents
$ws.Range("L10").Value2 = "@"
$ws.Range("M10").Value2 = "Free-air"
$ws.Range("S10").Value2 = "@"
$ws.Range("T10").Value2 = "Bouguer Simples"

# Clear the cells that are no longer part of any group (fully reset, not just clear contents)
$ws.Range("I10").Clear()
$ws.Range("J10").ClearContents()
$ws.Range("N10").Clear()
$ws.Range("O10").Clear()
$ws.Range("P10").Clear()
$ws.Range("U10").Clear()
$ws.Range("V10").Clear()
$ws.Range("W10").Clear()

# Row 10 is now a touch shorter
$ws.Rows.Item(10).RowHeight = 13.8

# --- Misc ---
$ws.Range("J11").Select()
